$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 16.34
$ws.Range("A8").Value = -22.188
$ws.Range("A10").Value = -21.856
$ws.Range("A12").Value = -21.578
$ws.Range("E14").Value = 16.982
$ws.Range("E15").Value = 16.361
$ws.Range("A18").Value = -22.095
$ws.Range("E18").Value = 16.386
$ws.Range("E20").Value = 16.455
$ws.Range("A25").Value = -21.808
$ws.Range("E29").Value = 17.05
$ws.Range("E30").Value = 16.182
$ws.Range("E31").Value = 16.225
$ws.Range("E35").Value = 16.589
$ws.Range("A37").Value = -20.232
$ws.Range("E40").Value = 16.627
$ws.Range("E44").Value = 16.574
$ws.Range("E50").Value = 16.319
$ws.Range("E54").Value = 16.72
$ws.Range("A55").Value = -22.166
$ws.Range("A68").Value = -21.603
$ws.Range("E68").Value = 17.333
$ws.Range("E76").Value = 16.558
$ws.Range("A77").Value = -20.48
$ws.Range("A78").Value = -20.008
$ws.Range("A79").Value = -21.57
$ws.Range("A80").Value = -20.193
$ws.Range("A81").Value = -21.782
$ws.Range("A82").Value = -22.152
$ws.Range("A84").Value = -22.047
$ws.Range("E87").Value = 16.366
$ws.Range("E88").Value = 16.278
$ws.Range("E92").Value = 17.914
$ws.Range("E96").Value = 16.325
$ws.Range("E98").Value = 16.293
$ws.Range("A101").Value = -21.295
$ws.Range("E101").Value = 16.936
$ws.Range("A102").Value = -20.501
$ws.Range("E102").Value = 16.459
